$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.120179
$ws.Range("H2").Value = 0.360537
$ws.Range("I2").Value = 0.04921086431616203
$ws.Range("J2").Value = 0.04921086431616202
$ws.Range("M2").Value = 0.7521946666666667
$ws.Range("N2").Value = 2.256584
$ws.Range("O2").Value = 0.07361670343069449
$ws.Range("P2").Value = 0.0736167034306945
$ws.Range("Q2").Value = 0.09039800284533334
$ws.Range("R2").Value = 0.8135820256080001
$ws.Range("S2").Value = 0.003622741603931046
$ws.Range("T2").Value = 0.003622741603931046

# Row 3
$ws.Range("G3").Value = 0.120179
$ws.Range("H3").Value = 0.360537
$ws.Range("I3").Value = 0.04921086431616203
$ws.Range("J3").Value = 0.04921086431616202
$ws.Range("O3").Value = 0.6908862423022597
$ws.Range("P3").Value = 0.6908862423022598
$ws.Range("Q3").Value = 0.8483772511796667
$ws.Range("R3").Value = 7.635395260617
$ws.Range("S3").Value = 0.03399910912783954
$ws.Range("T3").Value = 0.03399910912783954

# Row 4
$ws.Range("G4").Value = 0.120179
$ws.Range("H4").Value = 0.360537
$ws.Range("I4").Value = 0.04921086431616203
$ws.Range("J4").Value = 0.04921086431616202
$ws.Range("M4").Value = 2.406242333333334
$ws.Range("N4").Value = 7.218727
$ws.Range("O4").Value = 0.2354970542670457
$ws.Range("P4").Value = 0.2354970542670457
$ws.Range("Q4").Value = 0.2891797973776667
$ws.Range("R4").Value = 2.602618176399
$ws.Range("S4").Value = 0.01158901358439143
$ws.Range("T4").Value = 0.01158901358439143

# Row 5
$ws.Range("G5").Value = 1.522503666666667
$ws.Range("H5").Value = 4.567511000000001
$ws.Range("I5").Value = 0.6234343883806033
$ws.Range("J5").Value = 0.6234343883806033
$ws.Range("M5").Value = 0.7521946666666667
$ws.Range("N5").Value = 2.256584
$ws.Range("O5").Value = 0.07361670343069449
$ws.Range("P5").Value = 0.0736167034306945
$ws.Range("Q5").Value = 1.145219138047111
$ws.Range("R5").Value = 10.306972242424
$ws.Range("S5").Value = 0.04589518447791128
$ws.Range("T5").Value = 0.04589518447791129

# Row 6
$ws.Range("G6").Value = 1.522503666666667
$ws.Range("H6").Value = 4.567511000000001
$ws.Range("I6").Value = 0.6234343883806033
$ws.Range("J6").Value = 0.6234343883806033
$ws.Range("O6").Value = 0.6908862423022597
$ws.Range("P6").Value = 0.6908862423022598
$ws.Range("R6").Value = 96.73002172375101
$ws.Range("S6").Value = 0.4307222419102826
$ws.Range("T6").Value = 0.4307222419102826

# Row 7
$ws.Range("G7").Value = 1.522503666666667
$ws.Range("H7").Value = 4.567511000000001
$ws.Range("I7").Value = 0.6234343883806033
$ws.Range("J7").Value = 0.6234343883806033
$ws.Range("M7").Value = 2.406242333333334
$ws.Range("N7").Value = 7.218727
$ws.Range("O7").Value = 0.2354970542670457
$ws.Range("P7").Value = 0.2354970542670457
$ws.Range("Q7").Value = 3.663512775388556
$ws.Range("R7").Value = 32.97161497849701
$ws.Range("S7").Value = 0.1468169619924094
$ws.Range("T7").Value = 0.1468169619924094

# Row 8
$ws.Range("I8").Value = 0.3273547473032347
$ws.Range("J8").Value = 0.3273547473032347
$ws.Range("M8").Value = 0.7521946666666667
$ws.Range("N8").Value = 2.256584
$ws.Range("O8").Value = 0.07361670343069449
$ws.Range("P8").Value = 0.0736167034306945
$ws.Range("Q8").Value = 0.6013350057831111
$ws.Range("R8").Value = 5.412015052048
$ws.Range("S8").Value = 0.02409877734885216
$ws.Range("T8").Value = 0.02409877734885217

# Row 9
$ws.Range("I9").Value = 0.3273547473032347
$ws.Range("J9").Value = 0.3273547473032347
$ws.Range("O9").Value = 0.6908862423022597
$ws.Range("P9").Value = 0.6908862423022598
$ws.Range("S9").Value = 0.2261648912641376
$ws.Range("T9").Value = 0.2261648912641377

# Row 10
$ws.Range("I10").Value = 0.3273547473032347
$ws.Range("J10").Value = 0.3273547473032347
$ws.Range("M10").Value = 2.406242333333334
$ws.Range("N10").Value = 7.218727
$ws.Range("O10").Value = 0.2354970542670457
$ws.Range("P10").Value = 0.2354970542670457
$ws.Range("Q10").Value = 1.923647975121556
$ws.Range("S10").Value = 0.0770910786902449
$ws.Range("T10").Value = 0.0770910786902449
